$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty cells C5 and D5 with values
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 30

# Update the active cell / selection to D7 (was E7)
$ws.Range("D7").Select()
